$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header text updates (new report week / new volume number)
# ------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/27/2025  Through  2/2/2025"

# ------------------------------------------------------------------
# Cells that flip from a numeric style to a "no data" text style
# (style 13, shared text "0" / "***.*"). We copy a cell that already
# carries the desired text + style so the engine keeps the shared
# string reference and style index in sync (a direct .Value = "0"
# assignment gets auto-coerced back to the number 0).
# ------------------------------------------------------------------
$ws.Range("D15").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4104) | Out-Null
$ws.Range("D15").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null

$ws.Range("D15").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4104) | Out-Null
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null

$ws.Range("H15").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4104) | Out-Null
$ws.Range("H15").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null

$ws.Range("D15").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4104) | Out-Null
$ws.Range("D15").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null

# ------------------------------------------------------------------
# Cells that flip the other way: from the "no data" text placeholder
# back to a real number, so they need the ordinary numeric styles
# (14 / 15) restored before the new value is written.
# ------------------------------------------------------------------
$ws.Range("C15").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("D18").Value = 4

$ws.Range("H23").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = -50

$ws.Range("C15").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = 1

$ws.Range("H23").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = -100

# ------------------------------------------------------------------
# Plain numeric value updates (style/type unchanged)
# ------------------------------------------------------------------
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 3
$ws.Range("L15").Value = 50
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -33.333333333333
$ws.Range("J16").Value = 6
$ws.Range("K16").Value = -33.333333333333
$ws.Range("L16").Value = -50
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 18.181818181818
$ws.Range("I17").Value = 15
$ws.Range("J17").Value = 13
$ws.Range("K17").Value = 15.384615384615
$ws.Range("L17").Value = -37.5
$ws.Range("C18").Value = 2
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 9
$ws.Range("J18").Value = 8
$ws.Range("K18").Value = 12.5
$ws.Range("L18").Value = 50
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = -19.230769230769
$ws.Range("I19").Value = 26
$ws.Range("J19").Value = 33
$ws.Range("K19").Value = -21.212121212121
$ws.Range("L19").Value = -3.703703703703
$ws.Range("C20").Value = 3
$ws.Range("F20").Value = 12
$ws.Range("H20").Value = 300
$ws.Range("I20").Value = 12
$ws.Range("K20").Value = 300
$ws.Range("L20").Value = 100
$ws.Range("C21").Value = 13
$ws.Range("E21").Value = -18.75
$ws.Range("F21").Value = 61
$ws.Range("G21").Value = 53
$ws.Range("H21").Value = 15.094339622641
$ws.Range("I21").Value = 69
$ws.Range("J21").Value = 63
$ws.Range("K21").Value = 9.523809523809
$ws.Range("L21").Value = -5.479452054794
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -13.793103448275
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 113
$ws.Range("H24").Value = 7.079646017699
$ws.Range("I24").Value = 138
$ws.Range("J24").Value = 143
$ws.Range("K24").Value = -3.496503496503
$ws.Range("L24").Value = 18.965517241379
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 79
$ws.Range("G25").Value = 70
$ws.Range("H25").Value = 12.857142857142
$ws.Range("I25").Value = 92
$ws.Range("J25").Value = 86
$ws.Range("K25").Value = 6.976744186046
$ws.Range("L25").Value = 21.052631578947
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 32
$ws.Range("H26").Value = -15.789473684210
$ws.Range("I26").Value = 38
$ws.Range("J26").Value = 47
$ws.Range("K26").Value = -19.148936170212
$ws.Range("L26").Value = -24
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 3
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = -40
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = 33.333333333333
$ws.Range("L28").Value = -33.333333333333
$ws.Range("G31").Value = 1
